# "pr pa process O2" -- rename the GA sheet to "Report", restate the Report
# sheet's columns in terms of energy export / grid & plant availability /
# POA / PR, and rebuild the Raw Data sheet around the new PQM export +
# timestamp + POA layout (with extra 5-minute-granularity rows appended).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the "GA" sheet to "Report"
# ---------------------------------------------------------------------
$report = $wb.Worksheets.Item(1)
$report.Name = "Report"

# ---------------------------------------------------------------------
# 2. Report sheet: new headers + extra PR/POA columns
# ---------------------------------------------------------------------
$report.Range("A1").Value = "Timestamp"
$report.Range("B1").Value = "energyExport"
$report.Range("C1").Value = "Grid Availability (%)"
$report.Range("D1").Value = "Plant Availability (%)"
$report.Range("E1").Value = "POA(kWh/m2)"
$report.Range("F1").Value = "PR (%)"

# Row 2
$report.Range("D2").Value = -0.04166666666666667
$report.Range("E2").Value = 0
$report.Range("F2").Value = 0

# Row 3
$report.Range("D3").Value = -0.04166666666666667
$report.Range("E3").Value = 0
$report.Range("F3").Value = 0

# ---------------------------------------------------------------------
# 3. Raw Data sheet: rebuild with new headers/columns and rows
# ---------------------------------------------------------------------
$raw = $wb.Worksheets.Item(2)
$raw.Cells.ClearContents()

$raw.Range("A1").Value = "Actual Energy Export(kWh)"
$raw.Range("B1").Value = "PQM Active Total Export(kWh)"
$raw.Range("C1").Value = "Timestamp"
$raw.Range("D1").Value = "POA (w/m2)"

# A, B, C(, D) data rows. C holds date/timestamp text -- force text storage
# (leading apostrophe) so it isn't silently reinterpreted as a date serial,
# then reset the style back to Normal so no stray number format sticks.
$rows = @(
    @(2000, 2000, "2022-03-29", $null),
    @(0,    2000, "2022-03-29", $null),
    @(0,    2000, "2022-03-29", $null),
    @(0,    2000, "2022-03-29", $null),
    @(0,    2000, "2022-03-29", $null),
    @(2000, 2000, "2022-04-30", $null),
    @(0,    2000, "2022-04-30", $null),
    @(0,    2000, "2022-04-30", $null),
    @(0,    2000, "2022-04-30", $null),
    @(0,    2000, "2022-04-30", $null),
    @(0,    2000, "2022-03-29 00:55:00", 0),
    @(0,    2000, "2022-03-29 01:00:00", 0),
    @(0,    2000, "2022-03-29 01:20:00", 0),
    @(0,    2000, "2022-03-29 01:14:00", 0),
    @(0,    2000, "2022-03-29 00:57:00", 0),
    @(0,    2000, "2022-04-30 05:30:00", 0),
    @(0,    2000, "2022-04-30 06:20:00", 0),
    @(0,    2000, "2022-04-30 09:06:40", 0),
    @(0,    2000, "2022-04-30 06:43:20", 0),
    @(0,    2000, "2022-04-30 10:00:00", 0)
)

$r = 2
foreach ($row in $rows) {
    $raw.Cells.Item($r, 1).Value = $row[0]
    $raw.Cells.Item($r, 2).Value = $row[1]
    $raw.Cells.Item($r, 3).Value = "'" + $row[2]

    if ($null -ne $row[3]) {
        $raw.Cells.Item($r, 4).Value = $row[3]
    }

    $r = $r + 1
}

# The leading apostrophe above forces text storage (so the dates aren't
# reinterpreted as date serials); strip the resulting quote-prefix style
# back off in one pass so the cells end up plain "Normal"-styled text.
$raw.Range("C2:C21").Style = "Normal"
